$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get new values that look like plain numbers
# (e.g. "311.49"). Force Text format first so Excel keeps them as
# literal strings instead of silently converting them to numbers.
$textCells = @("D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D16", "D18", "D19", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '23.919.86'
$ws.Range("E2").Value = '  -2.09%  '
$ws.Range("D3").Value = '1.654.77'
$ws.Range("E3").Value = '  -0.86%  '
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '311.49'
$ws.Range("E5").Value = '  -0.28%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '0.3894'
$ws.Range("E7").Value = '  -1.90%  '
$ws.Range("D8").Value = '0.3822'
$ws.Range("E8").Value = '  -2.66%  '
$ws.Range("D9").Value = '51.68'
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").Value = '1.346'
$ws.Range("E10").Value = '  -3.47%  '
$ws.Range("D11").Value = '1.002'
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '0.08479'
$ws.Range("E12").Value = '  -1.01%  '
$ws.Range("D13").Value = '24.02'
$ws.Range("E13").Value = '  -2.09%  '
$ws.Range("D14").Value = '7.055'
$ws.Range("E14").Value = '  -3.36%  '
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").Value = '0.00001317'
$ws.Range("E16").Value = '  -1.39%  '
$ws.Range("D17").Value = '1.658.91'
$ws.Range("E17").Value = '  -0.29%  '
$ws.Range("D18").Value = '94.22'
$ws.Range("E18").Value = '  -0.79%  '
$ws.Range("D19").Value = '0.07001'
$ws.Range("E19").Value = '  -0.45%  '
$ws.Range("E20").Value = '  -5.03%  '
$ws.Range("E21").Value = '  -0.14%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '13.70'
$ws.Range("E23").Value = '  -0.46%  '
$ws.Range("D24").Value = '23.910.05'
$ws.Range("E24").Value = '  -2.17%  '
$ws.Range("D25").Value = '2.435'
$ws.Range("E25").Value = '  -1.94%  '
$ws.Range("D26").Value = '2.964'
$ws.Range("E26").Value = '  -3.37%  '
$ws.Range("D27").Value = '22.07'
$ws.Range("E27").Value = '  -2.10%  '
$ws.Range("D28").Value = '153.49'
$ws.Range("E28").Value = '  -2.43%  '
$ws.Range("D29").Value = '5.438'
$ws.Range("E29").Value = '  -0.35%  '
$ws.Range("D30").Value = '137.71'
$ws.Range("E30").Value = '  -3.43%  '
$ws.Range("D31").Value = '7.811'
$ws.Range("E31").Value = '  -2.79%  '
$ws.Range("D32").Value = '2.488'
$ws.Range("E32").Value = '  -2.32%  '
$ws.Range("D33").Value = '1.837.85'
$ws.Range("E33").Value = '  -0.41%  '
$ws.Range("D34").Value = '0.08150'
$ws.Range("E34").Value = '  -1.56%  '
$ws.Range("E35").Value = '  -5.57%  '
$ws.Range("D36").Value = '0.02922'
$ws.Range("E36").Value = '  -5.92%  '
$ws.Range("D37").Value = '6.660'
$ws.Range("E37").Value = '  -3.86%  '
$ws.Range("D38").Value = '10.80'
$ws.Range("E38").Value = '  -3.64%  '
$ws.Range("D39").Value = '0.2676'
$ws.Range("E39").Value = '  -3.22%  '
$ws.Range("D40").Value = '0.09128'
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("D41").Value = '13.62'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("D42").Value = '0.7571'
$ws.Range("E42").Value = '  -1.87%  '
$ws.Range("D43").Value = '1.428'
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("D44").Value = '16.47'
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("D45").Value = '0.6942'
$ws.Range("E45").Value = '  -2.53%  '
$ws.Range("D46").Value = '2.451'
$ws.Range("E46").Value = '  -3.83%  '
$ws.Range("D47").Value = '4.105'
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").Value = '0.9989'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("D49").Value = '0.08292'
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("D50").Value = '133.32'
$ws.Range("E50").Value = '  -2.69%  '
$ws.Range("E51").Value = '  -3.17%  '
